$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "59.151.84"
Set-TextValue "E2" "  -1.71%  "
Set-TextValue "D3" "2.563.03"
Set-TextValue "E3" "  -2.00%  "
Set-TextValue "E4" "  -0.03%  "
Set-TextValue "D5" "549.96"
Set-TextValue "E5" "  -3.23%  "
Set-TextValue "D6" "139.96"
Set-TextValue "E6" "  -3.65%  "
Set-TextValue "D7" "1.00"
Set-TextValue "E7" "  +0.29%  "
Set-TextValue "D8" "0.591"
Set-TextValue "E8" "  -1.58%  "
Set-TextValue "D9" "2.566.87"
Set-TextValue "E9" "  -2.62%  "
Set-TextValue "D10" "6.68"
Set-TextValue "E10" "  -1.10%  "
Set-TextValue "D11" "0.103"
Set-TextValue "E11" "  -0.99%  "
Set-TextValue "E12" "  +6.83%  "
Set-TextValue "D13" "0.352"
Set-TextValue "E13" "  +2.59%  "
Set-TextValue "D14" "3.021.17"
Set-TextValue "E14" "  -1.84%  "
Set-TextValue "D15" "59.216.58"
Set-TextValue "E15" "  -1.57%  "
Set-TextValue "D16" "22.98"
Set-TextValue "E16" "  +4.11%  "
Set-TextValue "D17" "0.0000136"
Set-TextValue "E17" "  -0.70%  "
Set-TextValue "D18" "2.579.60"
Set-TextValue "E18" "  -2.35%  "
Set-TextValue "D19" "4.53"
Set-TextValue "E19" "  +0.10%  "
Set-TextValue "D20" "336.33"
Set-TextValue "E20" "  -1.18%  "
Set-TextValue "D21" "10.24"
Set-TextValue "E21" "  -1.18%  "
Set-TextValue "D22" "6.39"
Set-TextValue "E22" "  +0.81%  "
Set-TextValue "D23" "0.996"
Set-TextValue "E23" "  -0.29%  "
Set-TextValue "D24" "0.474"
Set-TextValue "E24" "  +5.98%  "
Set-TextValue "D25" "62.63"
Set-TextValue "E25" "  -4.60%  "
Set-TextValue "D26" "1.00"
Set-TextValue "E26" "  +0.25%  "
Set-TextValue "D27" "0.158"
Set-TextValue "E27" "  -3.71%  "
Set-TextValue "D28" "7.37"
Set-TextValue "E28" "  +0.42%  "
Set-TextValue "D29" "0.0₃0761"
Set-TextValue "E29" "  -3.91%  "
Set-TextValue "E30" "  +0.07%  "
Set-TextValue "D31" "6.13"
Set-TextValue "E31" "  +0.10%  "
Set-TextValue "E32" "  -2.22%  "
Set-TextValue "D33" "158.48"
Set-TextValue "E33" "  -0.58%  "
Set-TextValue "D34" "18.97"
Set-TextValue "E34" "  -0.75%  "
Set-TextValue "D35" "4.07"
Set-TextValue "E35" "  -0.46%  "
Set-TextValue "D36" "1.16"
Set-TextValue "E36" "  +1.61%  "
Set-TextValue "D37" "0.890"
Set-TextValue "E37" "  +0.47%  "
Set-TextValue "D38" "37.40"
Set-TextValue "E38" "  -0.21%  "
Set-TextValue "D39" "0.845"
Set-TextValue "E39" "  -3.89%  "
Set-TextValue "D40" "1.46"
Set-TextValue "E40" "  -2.73%  "
Set-TextValue "D41" "3.65"
Set-TextValue "E41" "  +0.15%  "
Set-TextValue "D42" "285.18"
Set-TextValue "E42" "  -3.75%  "
Set-TextValue "D43" "135.79"
Set-TextValue "E43" "  +7.19%  "
Set-TextValue "E44" "  +0.35%  "
Set-TextValue "D45" "0.0966"
Set-TextValue "E45" "  -1.18%  "
Set-TextValue "B46" "WhiteBITCoin"
Set-TextValue "C46" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D46" "10.67"
Set-TextValue "E46" "  +0.01%  "
Set-TextValue "B47" "Mantle"
Set-TextValue "C47" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D47" "0.588"
Set-TextValue "E47" "  -2.30%  "
Set-TextValue "D48" "0.0527"
Set-TextValue "E48" "  -2.54%  "
Set-TextValue "D49" "0.0232"
Set-TextValue "E49" "  -0.91%  "
Set-TextValue "D50" "1.962.22"
Set-TextValue "E50" "  +0.18%  "
Set-TextValue "D51" "18.52"
Set-TextValue "E51" "  +0.03%  "
